$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Clear-Cell($ws, $addr) {
    $ws.Range($addr).ClearContents()
}

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
Set-Cell $ws "H32" 1367
Set-Cell $ws "I32" 1367
Set-Cell $ws "J32" 0
Set-Cell $ws "K32" 1367
Set-Cell $ws "L32" 0
Set-Cell $ws "M32" -1041
Clear-Cell $ws "N32"
Set-Cell $ws "H70" 2377.6
Set-Cell $ws "I70" 0
Set-Cell $ws "J70" 2377.6
Set-Cell $ws "K70" 0
Set-Cell $ws "L70" 7132.799999999999
Clear-Cell $ws "M70"
Set-Cell $ws "N70" -7672.799999999999
Set-Cell $ws "H73" 2377.6
Set-Cell $ws "I73" 0
Set-Cell $ws "J73" 2377.6
Set-Cell $ws "K73" 0
Set-Cell $ws "L73" 7132.799999999999
Clear-Cell $ws "M73"
Set-Cell $ws "N73" -9004.799999999999
Set-Cell $ws "H138" 2247.1953
Set-Cell $ws "I138" 1405.1708
Set-Cell $ws "J138" 2997.6956
Set-Cell $ws "K138" 4215.512400000001
Set-Cell $ws "L138" 8993.086800000001
Set-Cell $ws "M138" 924.4875999999995
Set-Cell $ws "N138" -19273.0868
Set-Cell $ws "H141" 1734.375
Set-Cell $ws "I141" 961.6667
Set-Cell $ws "J141" 4052.5
Set-Cell $ws "K141" 2885.0001
Set-Cell $ws "L141" 12157.5
Set-Cell $ws "M141" 2294.9999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
Set-Cell $ws "H122" 715.6
Set-Cell $ws "I122" 632.86206
Set-Cell $ws "J122" 1115.5
Set-Cell $ws "K122" 1898.58618
Set-Cell $ws "L122" 3346.5
Set-Cell $ws "M122" 551.4138199999998
Set-Cell $ws "N122" -8246.5
Set-Cell $ws "H128" 29333.334
Set-Cell $ws "I128" 0
Set-Cell $ws "J128" 29333.334
Set-Cell $ws "K128" 0
Set-Cell $ws "L128" 29333.334
Set-Cell $ws "N128" -39293.334
Set-Cell $ws "H132" 152930.38
Set-Cell $ws "I132" 23209.426
Set-Cell $ws "J132" 503176.94
Set-Cell $ws "K132" 69628.27799999999
Set-Cell $ws "L132" 1509530.82
Set-Cell $ws "M132" -67098.27799999999
Set-Cell $ws "N132" -1514590.82

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
Set-Cell $ws "H22" 1272.8334
Set-Cell $ws "I22" 4209.6
Set-Cell $ws "J22" 500
Set-Cell $ws "K22" 4209.6
Set-Cell $ws "L22" 500
Set-Cell $ws "M22" -4036.6
Set-Cell $ws "N22" -846
Set-Cell $ws "H80" 11408.294
Set-Cell $ws "I80" 21172.223
Set-Cell $ws "J80" 423.875
Set-Cell $ws "K80" 21172.223
Set-Cell $ws "L80" 423.875
Set-Cell $ws "M80" -20174.223
Set-Cell $ws "N80" -2419.875
Set-Cell $ws "H83" 11408.294
Set-Cell $ws "I83" 21172.223
Set-Cell $ws "J83" 423.875
Set-Cell $ws "K83" 105861.115
Set-Cell $ws "L83" 2119.375
Set-Cell $ws "M83" -100869.115
Set-Cell $ws "N83" -12103.375
Set-Cell $ws "H128" 1549
Set-Cell $ws "I128" 1549
Set-Cell $ws "J128" 0
Set-Cell $ws "K128" 4647
Set-Cell $ws "L128" 0
Set-Cell $ws "M128" -2157
Set-Cell $ws "H134" 41710944
Set-Cell $ws "I134" 2829.647
Set-Cell $ws "J134" 143002080
Set-Cell $ws "K134" 8488.940999999999
Set-Cell $ws "L134" 429006240
Set-Cell $ws "M134" -5953.940999999999
Set-Cell $ws "N134" -429011310

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
Set-Cell $ws "H31" 4209.7256
Set-Cell $ws "I31" 1278.9286
Set-Cell $ws "J31" 10364.4
Set-Cell $ws "K31" 1278.9286
Set-Cell $ws "L31" 10364.4
Set-Cell $ws "M31" -983.9286
Set-Cell $ws "N31" -10954.4
Set-Cell $ws "H34" 4209.7256
Set-Cell $ws "I34" 1278.9286
Set-Cell $ws "J34" 10364.4
Set-Cell $ws "K34" 1278.9286
Set-Cell $ws "L34" 10364.4
Set-Cell $ws "M34" -1076.9286
Set-Cell $ws "N34" -10768.4
Set-Cell $ws "H62" 2502.8572
Set-Cell $ws "I62" 2502.8572
Set-Cell $ws "J62" 0
Set-Cell $ws "K62" 2502.8572
Set-Cell $ws "L62" 0
Set-Cell $ws "M62" -1878.8572
Clear-Cell $ws "N62"
Set-Cell $ws "H65" 2502.8572
Set-Cell $ws "I65" 2502.8572
Set-Cell $ws "J65" 0
Set-Cell $ws "K65" 12514.286
Set-Cell $ws "L65" 0
Set-Cell $ws "M65" -9394.286
Clear-Cell $ws "N65"
Set-Cell $ws "H86" 111128664
Set-Cell $ws "I86" 111128664
Set-Cell $ws "J86" 0
Set-Cell $ws "K86" 111128664
Set-Cell $ws "L86" 0
Set-Cell $ws "M86" -111127541
Set-Cell $ws "H89" 111128664
Set-Cell $ws "I89" 111128664
Set-Cell $ws "J89" 0
Set-Cell $ws "K89" 555643320
Set-Cell $ws "L89" 0
Set-Cell $ws "M89" -555637704
Set-Cell $ws "H132" 23818.613
Set-Cell $ws "I132" 37960.703
Set-Cell $ws "J132" 1357.6471
Set-Cell $ws "K132" 113882.109
Set-Cell $ws "L132" 4072.9413
Set-Cell $ws "M132" -111352.109
Set-Cell $ws "N132" -9132.941299999999
Set-Cell $ws "H133" 42000
Set-Cell $ws "I133" 0
Set-Cell $ws "J133" 42000
Set-Cell $ws "K133" 0
Set-Cell $ws "L133" 42000
Set-Cell $ws "N133" -47060

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
Set-Cell $ws "H126" 1857.6923
Set-Cell $ws "I126" 1325
Set-Cell $ws "J126" 1954.5454
Set-Cell $ws "K126" 3975
Set-Cell $ws "L126" 5863.6362
Set-Cell $ws "M126" -1505
Set-Cell $ws "N126" -10803.6362

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
Set-Cell $ws "H136" 910839.4399999999
Set-Cell $ws "I136" 2001256.8
Set-Cell $ws "J136" 2158.3333
Set-Cell $ws "K136" 6003770.4
Set-Cell $ws "L136" 6474.999899999999
Set-Cell $ws "M136" -6001220.4
Set-Cell $ws "N136" -11574.9999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
Set-Cell $ws "H11" 75003.75
Set-Cell $ws "I11" 20000
Set-Cell $ws "J11" 93338.336
Set-Cell $ws "K11" 20000
Set-Cell $ws "L11" 93338.336
Set-Cell $ws "M11" -19858
Set-Cell $ws "N11" -93622.336
Set-Cell $ws "H54" 20077
Set-Cell $ws "I54" 0
Set-Cell $ws "J54" 20077
Set-Cell $ws "K54" 0
Set-Cell $ws "L54" 20077
Clear-Cell $ws "M54"
Set-Cell $ws "N54" -21117
Set-Cell $ws "H62" 5150
Set-Cell $ws "I62" 4866.6665
Set-Cell $ws "J62" 5215.385
Set-Cell $ws "K62" 4866.6665
Set-Cell $ws "L62" 5215.385
Set-Cell $ws "M62" -4242.6665
Set-Cell $ws "N62" -6463.385
Set-Cell $ws "H65" 5150
Set-Cell $ws "I65" 4866.6665
Set-Cell $ws "J65" 5215.385
Set-Cell $ws "K65" 24333.3325
Set-Cell $ws "L65" 26076.925
Set-Cell $ws "M65" -21213.3325
Set-Cell $ws "N65" -32316.925
Set-Cell $ws "H81" 2312.158
Set-Cell $ws "I81" 2766.3333
Set-Cell $ws "J81" 1903.4
Set-Cell $ws "K81" 5532.6666
Set-Cell $ws "L81" 3806.8
Set-Cell $ws "M81" -4471.6666
Set-Cell $ws "N81" -5928.8
Set-Cell $ws "H84" 2312.158
Set-Cell $ws "I84" 2766.3333
Set-Cell $ws "J84" 1903.4
Set-Cell $ws "K84" 27663.333
Set-Cell $ws "L84" 19034
Set-Cell $ws "M84" -22359.333
Set-Cell $ws "N84" -29642
